$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43; this shifts existing rows 43-185 down to 44-186
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with its data
$ws.Cells.Item(43, 1).Value = 10
$ws.Cells.Item(43, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value = "La Araucanía"
$ws.Cells.Item(43, 4).Value = 44624
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 6).Value = 100112005
$ws.Cells.Item(43, 7).Value = "Puerro"
$ws.Cells.Item(43, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 10
$ws.Cells.Item(43, 11).Value = 12000
$ws.Cells.Item(43, 12).Value = 12000
$ws.Cells.Item(43, 13).Value = 12000
$ws.Cells.Item(43, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(43, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(43, 16).Value = 1000
$ws.Cells.Item(43, 17).Value = 12
$ws.Cells.Item(43, 18).Value = "Hortaliza"
